$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.201.67"
$ws.Range("E2").Value = "  +0.62%  "

$ws.Range("D3").Value = "1.836.09"
$ws.Range("E3").Value = "  +0.37%  "

$ws.Range("D4").Value = "0.9985"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "242.18"
$ws.Range("E5").Value = "  +1.38%  "

$ws.Range("D6").Value = "0.6606"
$ws.Range("E6").Value = "  -0.47%  "

$ws.Range("D7").Value = "0.9995"
$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").Value = "0.07443"
$ws.Range("E8").Value = "  +0.96%  "

$ws.Range("D9").Value = "0.2931"
$ws.Range("E9").Value = "  -0.50%  "

$ws.Range("D10").Value = "23.00"
$ws.Range("E10").Value = "  +1.47%  "

$ws.Range("D11").Value = "0.07755"
$ws.Range("E11").Value = "  +1.55%  "

$ws.Range("D12").Value = "1.817.32"
$ws.Range("E12").Value = "  -0.97%  "

$ws.Range("D13").Value = "4.988"
$ws.Range("E13").Value = "  -0.30%  "

$ws.Range("D14").Value = "0.6667"
$ws.Range("E14").Value = "  -0.64%  "

$ws.Range("D15").Value = "82.90"
$ws.Range("E15").Value = "  -3.68%  "

$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "2.336.68"
$ws.Range("E16").Value = "  +12.33%  "

$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").Value = "6.129"
$ws.Range("E17").Value = "  +0.25%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.000008589"
$ws.Range("E18").Value = "  +4.80%  "

$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "29.194.32"
$ws.Range("E19").Value = "  +0.53%  "

$ws.Range("D20").Value = "227.29"
$ws.Range("E20").Value = "  +0.06%  "

$ws.Range("D21").Value = "12.49"
$ws.Range("E21").Value = "  +0.53%  "

$ws.Range("D22").Value = "0.9994"
$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("D23").Value = "7.139"
$ws.Range("E23").Value = "  -1.66%  "

$ws.Range("D24").Value = "0.9992"
$ws.Range("E24").Value = "  -0.12%  "

$ws.Range("D25").Value = "158.86"
$ws.Range("E25").Value = "  -0.87%  "

$ws.Range("D26").Value = "8.613"
$ws.Range("E26").Value = "  -0.37%  "

$ws.Range("D27").Value = "0.1399"
$ws.Range("E27").Value = "  -1.49%  "

$ws.Range("D28").Value = "17.96"
$ws.Range("E28").Value = "  +0.20%  "

$ws.Range("D29").Value = "1.515"
$ws.Range("E29").Value = "  +0.88%  "

$ws.Range("D30").Value = "4.121"
$ws.Range("E30").Value = "  -2.45%  "

$ws.Range("D31").Value = "4.052"
$ws.Range("E31").Value = "  -1.43%  "

$ws.Range("D32").Value = "1.195"
$ws.Range("E32").Value = "  -0.07%  "

$ws.Range("D33").Value = "0.05256"
$ws.Range("E33").Value = "  -1.79%  "

$ws.Range("D34").Value = "1.868"
$ws.Range("E34").Value = "  +1.48%  "

$ws.Range("D35").Value = "0.7370"
$ws.Range("E35").Value = "  -1.33%  "

$ws.Range("D36").Value = "1.147"
$ws.Range("E36").Value = "  +1.96%  "

$ws.Range("D37").Value = "2.654"
$ws.Range("E37").Value = "  -1.01%  "

$ws.Range("D38").Value = "1.309.51"
$ws.Range("E38").Value = "  +1.17%  "

$ws.Range("D39").Value = "0.01798"
$ws.Range("E39").Value = "  -0.32%  "

$ws.Range("D40").Value = "2.737"
$ws.Range("E40").Value = "  +1.25%  "

$ws.Range("D41").Value = "0.9239"
$ws.Range("E41").Value = "  +0.32%  "

$ws.Range("D42").Value = "0.08885"
$ws.Range("E42").Value = "  +13.14%  "

$ws.Range("D43").Value = "6.060"
$ws.Range("E43").Value = "  +0.70%  "

$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "2.135.27"
$ws.Range("E44").Value = "  +7.91%  "

$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "0.9989"
$ws.Range("E45").Value = "  +0.02%  "

$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "102.62"
$ws.Range("E46").Value = "  -1.19%  "

$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.00000000120"
$ws.Range("E47").Value = "  -2.38%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "0.5146"
$ws.Range("E48").Value = "  -0.60%  "

$ws.Range("D49").Value = "63.76"
$ws.Range("E49").Value = "  +0.74%  "

$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "1.758"
$ws.Range("E50").Value = "  +0.65%  "

$ws.Range("D51").Value = "0.05844"
$ws.Range("E51").Value = "  -1.13%  "
